$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 9.036864
$ws.Range("H2").Value = 27.110592
$ws.Range("I2").Value = 0.9206342953678062
$ws.Range("J2").Value = 0.9206342953678062
$ws.Range("M2").Value = 1.819857
$ws.Range("N2").Value = 5.459571
$ws.Range("O2").Value = 0.01485317462584607
$ws.Range("P2").Value = 0.01485317462584607
$ws.Range("Q2").Value = 16.445800208448
$ws.Range("R2").Value = 148.012201876032
$ws.Range("S2").Value = 0.01367434195564077
$ws.Range("T2").Value = 0.01367434195564077

# Row 3
$ws.Range("G3").Value = 9.036864
$ws.Range("H3").Value = 27.110592
$ws.Range("I3").Value = 0.9206342953678062
$ws.Range("J3").Value = 0.9206342953678062
$ws.Range("O3").Value = 0.726618572334523
$ws.Range("P3").Value = 0.7266185723345231
$ws.Range("Q3").Value = 804.5299519718399
$ws.Range("R3").Value = 7240.769567746559
$ws.Range("S3").Value = 0.6689499773423548
$ws.Range("T3").Value = 0.668949977342355

# Row 4
$ws.Range("G4").Value = 9.036864
$ws.Range("H4").Value = 27.110592
$ws.Range("I4").Value = 0.9206342953678062
$ws.Range("J4").Value = 0.9206342953678062
$ws.Range("M4").Value = 31.52924033333333
$ws.Range("N4").Value = 94.58772099999999
$ws.Range("O4").Value = 0.257333028084772
$ws.Range("P4").Value = 0.257333028084772
$ws.Range("Q4").Value = 284.9254569156479
$ws.Range("R4").Value = 2564.329112240831
$ws.Range("S4").Value = 0.236909610985688
$ws.Range("T4").Value = 0.236909610985688

# Row 5
$ws.Range("G5").Value = 9.036864
$ws.Range("H5").Value = 27.110592
$ws.Range("I5").Value = 0.9206342953678062
$ws.Range("J5").Value = 0.9206342953678062
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1464426666666667
$ws.Range("N5").Value = 0.439328
$ws.Range("O5").Value = 0.001195224954858853
$ws.Range("P5").Value = 0.001195224954858853
$ws.Range("Q5").Value = 1.323382462464
$ws.Range("R5").Value = 11.910442162176
$ws.Range("S5").Value = 0.001100365084122498
$ws.Range("T5").Value = 0.001100365084122498

# Row 6
$ws.Range("H6").Value = 0.730256
$ws.Range("I6").Value = 0.02479837836068326
$ws.Range("J6").Value = 0.02479837836068326
$ws.Range("M6").Value = 1.819857
$ws.Range("N6").Value = 5.459571
$ws.Range("O6").Value = 0.01485317462584607
$ws.Range("P6").Value = 0.01485317462584607
$ws.Range("Q6").Value = 0.442987164464
$ws.Range("R6").Value = 3.986884480176
$ws.Range("S6").Value = 0.0003683346442290307
$ws.Range("T6").Value = 0.0003683346442290309

# Row 7
$ws.Range("H7").Value = 0.730256
$ws.Range("I7").Value = 0.02479837836068326
$ws.Range("J7").Value = 0.02479837836068326
$ws.Range("O7").Value = 0.726618572334523
$ws.Range("P7").Value = 0.7266185723345231
$ws.Range("S7").Value = 0.018018962280651
$ws.Range("T7").Value = 0.01801896228065101

# Row 8
$ws.Range("H8").Value = 0.730256
$ws.Range("I8").Value = 0.02479837836068326
$ws.Range("J8").Value = 0.02479837836068326
$ws.Range("M8").Value = 31.52924033333333
$ws.Range("N8").Value = 94.58772099999999
$ws.Range("O8").Value = 0.257333028084772
$ws.Range("P8").Value = 0.257333028084772
$ws.Range("Q8").Value = 7.674805642952888
$ws.Range("R8").Value = 69.073250786576
$ws.Range("S8").Value = 0.006381441795146508
$ws.Range("T8").Value = 0.006381441795146509

# Row 9
$ws.Range("H9").Value = 0.730256
$ws.Range("I9").Value = 0.02479837836068326
$ws.Range("J9").Value = 0.02479837836068326
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1464426666666667
$ws.Range("N9").Value = 0.439328
$ws.Range("O9").Value = 0.001195224954858853
$ws.Range("P9").Value = 0.001195224954858853
$ws.Range("Q9").Value = 0.03564687866311111
$ws.Range("R9").Value = 0.320821907968
$ws.Range("S9").Value = 0.00002963964065672039
$ws.Range("T9").Value = 0.0000296396406567204

# Row 10
$ws.Range("G10").Value = 0.5126043333333333
$ws.Range("H10").Value = 1.537813
$ws.Range("I10").Value = 0.05222178061115199
$ws.Range("J10").Value = 0.05222178061115199
$ws.Range("M10").Value = 1.819857
$ws.Range("N10").Value = 5.459571
$ws.Range("O10").Value = 0.01485317462584607
$ws.Range("P10").Value = 0.01485317462584607
$ws.Range("Q10").Value = 0.932866584247
$ws.Range("R10").Value = 8.395799258223001
$ws.Range("S10").Value = 0.0007756592266900627
$ws.Range("T10").Value = 0.0007756592266900629

# Row 11
$ws.Range("G11").Value = 0.5126043333333333
$ws.Range("H11").Value = 1.537813
$ws.Range("I11").Value = 0.05222178061115199
$ws.Range("J11").Value = 0.05222178061115199
$ws.Range("O11").Value = 0.726618572334523
$ws.Range("P11").Value = 0.7266185723345231
$ws.Range("Q11").Value = 45.63591304209333
$ws.Range("R11").Value = 410.72321737884
$ws.Range("S11").Value = 0.03794531567244193
$ws.Range("T11").Value = 0.03794531567244194

# Row 12
$ws.Range("G12").Value = 0.5126043333333333
$ws.Range("H12").Value = 1.537813
$ws.Range("I12").Value = 0.05222178061115199
$ws.Range("J12").Value = 0.05222178061115199
$ws.Range("M12").Value = 31.52924033333333
$ws.Range("N12").Value = 94.58772099999999
$ws.Range("O12").Value = 0.257333028084772
$ws.Range("P12").Value = 0.257333028084772
$ws.Range("Q12").Value = 16.16202522157478
$ws.Range("R12").Value = 145.458226994173
$ws.Range("S12").Value = 0.01343838893664638
$ws.Range("T12").Value = 0.01343838893664638

# Row 13
$ws.Range("G13").Value = 0.5126043333333333
$ws.Range("H13").Value = 1.537813
$ws.Range("I13").Value = 0.05222178061115199
$ws.Range("J13").Value = 0.05222178061115199
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1464426666666667
$ws.Range("N13").Value = 0.439328
$ws.Range("O13").Value = 0.001195224954858853
$ws.Range("P13").Value = 0.001195224954858853
$ws.Range("Q13").Value = 0.07506714551822222
$ws.Range("R13").Value = 0.6756043096640001
$ws.Range("S13").Value = 0.00006241677537361303
$ws.Range("T13").Value = 0.00006241677537361305

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.02302366666666667
$ws.Range("H14").Value = 0.06907100000000001
$ws.Range("I14").Value = 0.002345545660358495
$ws.Range("J14").Value = 0.002345545660358496
$ws.Range("M14").Value = 1.819857
$ws.Range("N14").Value = 5.459571
$ws.Range("O14").Value = 0.01485317462584607
$ws.Range("P14").Value = 0.01485317462584607
$ws.Range("Q14").Value = 0.041899780949
$ws.Range("R14").Value = 0.3770980285410001
$ws.Range("S14").Value = 0.00003483879928620016
$ws.Range("T14").Value = 0.00003483879928620017

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.02302366666666667
$ws.Range("H15").Value = 0.06907100000000001
$ws.Range("I15").Value = 0.002345545660358495
$ws.Range("J15").Value = 0.002345545660358496
$ws.Range("O15").Value = 0.726618572334523
$ws.Range("P15").Value = 0.7266185723345231
$ws.Range("Q15").Value = 2.049740865586667
$ws.Range("R15").Value = 18.44766779028
$ws.Range("S15").Value = 0.001704317039075126
$ws.Range("T15").Value = 0.001704317039075126

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.02302366666666667
$ws.Range("H16").Value = 0.06907100000000001
$ws.Range("I16").Value = 0.002345545660358495
$ws.Range("J16").Value = 0.002345545660358496
$ws.Range("M16").Value = 31.52924033333333
$ws.Range("N16").Value = 94.58772099999999
$ws.Range("O16").Value = 0.257333028084772
$ws.Range("P16").Value = 0.257333028084772
$ws.Range("Q16").Value = 0.7259187196878889
$ws.Range("R16").Value = 6.533268477191
$ws.Range("S16").Value = 0.0006035863672911478
$ws.Range("T16").Value = 0.000603586367291148

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.02302366666666667
$ws.Range("H17").Value = 0.06907100000000001
$ws.Range("I17").Value = 0.002345545660358495
$ws.Range("J17").Value = 0.002345545660358496
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1464426666666667
$ws.Range("N17").Value = 0.439328
$ws.Range("O17").Value = 0.001195224954858853
$ws.Range("P17").Value = 0.001195224954858853
$ws.Range("Q17").Value = 0.003371647143111111
$ws.Range("R17").Value = 0.030344824288
$ws.Range("S17").Value = 0.00000280345470602136
$ws.Range("T17").Value = 0.000002803454706021361

